$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: clear F1:J1 header text but keep style
$ws.Range("F1:J1").ClearContents()

# Row 2: D2/E2 get new numeric values with #,##0 format; H2:J2 cleared but keep style
$ws.Range("D2").Value = 6192716
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("E2").Value = 6234953
$ws.Range("E2").NumberFormat = "#,##0"
$ws.Range("H2:J2").ClearContents()

# Row 3: H3,I3 removed; J3 cleared (keeps its style)
$ws.Range("H3:J3").ClearContents()

# Row 5: H5,I5 removed; J5 cleared (keeps its style)
$ws.Range("H5:J5").ClearContents()

# Row 7: A7 changes from number 14 to text "Chromosome inversion 14"
$ws.Range("A7").Value = "Chromosome inversion 14"
$ws.Range("H7:J7").ClearContents()

# Update selection
$ws.Range("G23").Select() | Out-Null
